{"js": "// Update formatting of the \"Table Caption\" style: remove the paragraph's\n// explicit \"space after\" override (0pt) so that the paragraph spacing goes\n// back to matching the 6pt / 120-twip value inherited from the base\n// \"Caption\" style.\nconst style = context.document.getStyles().getByNameOrNullObject(\"Table Caption\");\nstyle.load(\"isNullObject,nameLocal\");\nawait context.sync();\n\nif (!style.isNullObject) {\n  // 120 twips == 6 points, the \"space after\" value defined on the base\n  // \"Caption\" style that \"Table Caption\" previously overrode with 0.\n  style.paragraphFormat.spaceAfter = 6;\n  await context.sync();\n}\n", "ps1": "# Update formatting of the \"Table Caption\" style: remove the paragraph's\n# explicit \"space after\" override (0pt) so that paragraph spacing goes back\n# to matching the 6pt / 120-twip value inherited from the base \"Caption\"\n# style.\n$d = $word.ActiveDocument\n$style = $d.Styles(\"Table Caption\")\n\n# 120 twips == 6 points, the \"space after\" value defined on the base\n# \"Caption\" style that \"Table Caption\" previously overrode with 0.\n$style.ParagraphFormat.SpaceAfter = 6\n"}
